{"js": "// Changed 2015 table labels\n// - Remove the \"Abstract Title\" paragraph style\n// - Change the \"Abstract\" style's space-before from 5pt to 15pt (100 -> 300 twips)\n// - Remove the \"Footnote Block Text\" paragraph style\n//\n// Style handles are resolved by their position in the style collection, so\n// deleting one style shifts the indices of every style that follows it.\n// To avoid acting on a stale position, each style is looked up and used one\n// at a time, syncing in between so no handle is held across a delete of an\n// earlier-indexed style.\n\nconst styles = context.document.getStyles();\n\nconst abstractTitleStyle = styles.getByNameOrNullObject(\"Abstract Title\");\nabstractTitleStyle.load(\"isNullObject\");\nawait context.sync();\nif (!abstractTitleStyle.isNullObject) {\n  abstractTitleStyle.delete();\n  await context.sync();\n}\n\nconst abstractStyle = styles.getByNameOrNullObject(\"Abstract\");\nabstractStyle.load(\"isNullObject\");\nawait context.sync();\nif (!abstractStyle.isNullObject) {\n  abstractStyle.paragraphFormat.spaceBefore = 15;\n  await context.sync();\n}\n\nconst footnoteBlockTextStyle = styles.getByNameOrNullObject(\"Footnote Block Text\");\nfootnoteBlockTextStyle.load(\"isNullObject\");\nawait context.sync();\nif (!footnoteBlockTextStyle.isNullObject) {\n  footnoteBlockTextStyle.delete();\n  await context.sync();\n}\n", "ps1": "# Changed 2015 table labels\n# - Remove the \"Abstract Title\" paragraph style\n# - Change the \"Abstract\" style's space-before from 5pt (100 twips) to 15pt (300 twips)\n# - Remove the \"Footnote Block Text\" paragraph style\n\n$d = $word.ActiveDocument\n\n# Remove the \"Abstract Title\" paragraph style entirely\n$d.Styles(\"Abstract Title\").Delete()\n\n# Change the \"Abstract\" style's space-before from 5pt (100 twips) to 15pt (300 twips)\n$d.Styles(\"Abstract\").ParagraphFormat.SpaceBefore = 15\n\n# Remove the \"Footnote Block Text\" paragraph style entirely\n$d.Styles(\"Footnote Block Text\").Delete()\n"}
